$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 55, shifting the existing rows 55-79 down to 58-82.
$ws.Rows("55:57").Insert()

# New row 55: Maracuya, Especial, week of 44603
$ws.Cells.Item(55, 1).Value = 1
$ws.Cells.Item(55, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(55, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(55, 4).Value = 44603
$ws.Cells.Item(55, 5).Value = 15
$ws.Cells.Item(55, 6).Value = "Fruta"
$ws.Cells.Item(55, 7).Value = 100108
$ws.Cells.Item(55, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(55, 9).Value = 100108003
$ws.Cells.Item(55, 10).Value = "Maracuyá"
$ws.Cells.Item(55, 11).Value = "Sin especificar"
$ws.Cells.Item(55, 12).Value = "Especial"
$ws.Cells.Item(55, 13).Value = 70
$ws.Cells.Item(55, 14).Value = 44000
$ws.Cells.Item(55, 15).Value = 45000
$ws.Cells.Item(55, 16).Value = 44500
$ws.Cells.Item(55, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(55, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(55, 19).Value = 2225
$ws.Cells.Item(55, 20).Value = 20

# New row 56: Maracuya, Primera, week of 44603
$ws.Cells.Item(56, 1).Value = 1
$ws.Cells.Item(56, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(56, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(56, 4).Value = 44603
$ws.Cells.Item(56, 5).Value = 15
$ws.Cells.Item(56, 6).Value = "Fruta"
$ws.Cells.Item(56, 7).Value = 100108
$ws.Cells.Item(56, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(56, 9).Value = 100108003
$ws.Cells.Item(56, 10).Value = "Maracuyá"
$ws.Cells.Item(56, 11).Value = "Sin especificar"
$ws.Cells.Item(56, 12).Value = "Primera"
$ws.Cells.Item(56, 13).Value = 120
$ws.Cells.Item(56, 14).Value = 40000
$ws.Cells.Item(56, 15).Value = 41000
$ws.Cells.Item(56, 16).Value = 40500
$ws.Cells.Item(56, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(56, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(56, 19).Value = 2025
$ws.Cells.Item(56, 20).Value = 20

# New row 57: Maracuya, Segunda, week of 44603
$ws.Cells.Item(57, 1).Value = 1
$ws.Cells.Item(57, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(57, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(57, 4).Value = 44603
$ws.Cells.Item(57, 5).Value = 15
$ws.Cells.Item(57, 6).Value = "Fruta"
$ws.Cells.Item(57, 7).Value = 100108
$ws.Cells.Item(57, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(57, 9).Value = 100108003
$ws.Cells.Item(57, 10).Value = "Maracuyá"
$ws.Cells.Item(57, 11).Value = "Sin especificar"
$ws.Cells.Item(57, 12).Value = "Segunda"
$ws.Cells.Item(57, 13).Value = 120
$ws.Cells.Item(57, 14).Value = 37000
$ws.Cells.Item(57, 15).Value = 38000
$ws.Cells.Item(57, 16).Value = 37500
$ws.Cells.Item(57, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(57, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(57, 19).Value = 1875
$ws.Cells.Item(57, 20).Value = 20

# Apply the date-style (numFmtId 165) to the D cells of the new rows, matching
# the formatting used throughout the rest of the Fecha column.
$ws.Range("D55:D57").NumberFormat = $ws.Range("D58").NumberFormat
